$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "health" row (row 2) entirely; rows below shift up.
$ws.Rows.Item(2).Delete()
